$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns remain text, matching the
# original inline-string cell type, so numeric-looking values like
# "289.30" or "1.003" are not silently converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '21.571.65'
$ws.Range("E2").Value = '  -2.45%  '
$ws.Range("D3").Value = '1.531.89'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = '289.30'
$ws.Range("E6").Value = '  -0.84%  '
$ws.Range("D7").Value = '0.3894'
$ws.Range("E7").Value = '  -1.46%  '
$ws.Range("D8").Value = '0.3168'
$ws.Range("E8").Value = '  -2.31%  '
$ws.Range("D9").Value = '42.81'
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("D10").Value = '0.07155'
$ws.Range("E10").Value = '  -2.52%  '
$ws.Range("D11").Value = '1.068'
$ws.Range("E11").Value = '  -2.18%  '
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '5.735'
$ws.Range("E13").Value = '  +1.12%  '
$ws.Range("D14").Value = '18.19'
$ws.Range("E14").Value = '  -5.08%  '
$ws.Range("D15").Value = '6.570'
$ws.Range("E15").Value = '  -1.57%  '
$ws.Range("D16").Value = '1.537.77'
$ws.Range("E16").Value = '  -1.15%  '
$ws.Range("D17").Value = '0.00001085'
$ws.Range("E17").Value = '  -5.56%  '
$ws.Range("D18").Value = '0.06616'
$ws.Range("E18").Value = '  +0.22%  '
$ws.Range("D19").Value = '83.92'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").Value = '6.105'
$ws.Range("E21").Value = '  -3.42%  '
$ws.Range("D22").Value = '15.36'
$ws.Range("E22").Value = '  -3.15%  '
$ws.Range("E23").Value = '  -5.34%  '
$ws.Range("D24").Value = '2.380'
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("D25").Value = '21.566.22'
$ws.Range("E25").Value = '  -2.53%  '
$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").Value = '2.354'
$ws.Range("E26").Value = '  -3.75%  '
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '150.79'
$ws.Range("E27").Value = '  +1.68%  '
$ws.Range("D28").Value = '18.40'
$ws.Range("E28").Value = '  -1.65%  '
$ws.Range("D29").Value = '4.885'
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").Value = '1.709.12'
$ws.Range("E30").Value = '  -1.17%  '
$ws.Range("D31").Value = '116.76'
$ws.Range("E31").Value = '  -2.10%  '
$ws.Range("D32").Value = '6.012'
$ws.Range("E32").Value = '  +4.84%  '
$ws.Range("D33").Value = '0.9490'
$ws.Range("E33").Value = '  -7.63%  '
$ws.Range("D34").Value = '0.08020'
$ws.Range("E34").Value = '  -4.12%  '
$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").Value = '8.485'
$ws.Range("E35").Value = '  -6.73%  '
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '5.143'
$ws.Range("E36").Value = '  -0.54%  '
$ws.Range("D37").Value = '1.490'
$ws.Range("E37").Value = '  -8.66%  '
$ws.Range("D38").Value = '0.02203'
$ws.Range("E38").Value = '  -3.38%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05905'
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").Value = '11.22'
$ws.Range("E40").Value = '  +3.65%  '
$ws.Range("D41").Value = '0.2027'
$ws.Range("E41").Value = '  -2.07%  '
$ws.Range("D42").Value = '1.177'
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("D43").Value = '1.003'
$ws.Range("E43").Value = '  +0.37%  '
$ws.Range("D44").Value = '0.5753'
$ws.Range("E44").Value = '  -2.31%  '
$ws.Range("D45").Value = '13.07'
$ws.Range("E45").Value = '  -0.02%  '
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("D47").Value = '0.5530'
$ws.Range("E47").Value = '  -2.03%  '
$ws.Range("B48").Value = 'EOS'
$ws.Range("C48").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D48").Value = '1.161'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.884'
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("D50").Value = '115.69'
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("D51").Value = '0.06689'
$ws.Range("E51").Value = '  -2.67%  '
